$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.848.10"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "2.637.13"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("E4").Value = "  +0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "592.99"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "155.29"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.20%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("E9").Value = "  +6.08%  "
$ws.Range("E10").Value = "  +3.57%  "
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("E12").Value = "  +1.79%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "28.92"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +5.03%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000186"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +19.32%  "
$ws.Range("D15").Value = "3.110.92"
$ws.Range("E15").Value = "  +2.24%  "
$ws.Range("D16").Value = "64.771.76"
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("D17").Value = "2.645.70"
$ws.Range("E17").Value = "  +2.19%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "12.54"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +3.01%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.79"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.19%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "350.99"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.41%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.25"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +5.88%  "
$ws.Range("E22").Value = "  +0.06%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "68.04"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("E24").Value = "  +0.43%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.47"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.67%  "
$ws.Range("E26").Value = "  -2.23%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.10"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.33%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "0.0₃0934"
$ws.Range("E30").Value = "  +8.39%  "
$ws.Range("E31").Value = "  +2.59%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "509.33"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -8.00%  "
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("E34").Value = "  +6.86%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "6.21"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("E36").Value = "  +2.46%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "164.99"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.04%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "20.11"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.96%  "
$ws.Range("E39").Value = "  +4.47%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +0.07%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "42.23"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +6.44%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "164.36"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("E44").Value = "  +2.58%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0611"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +4.33%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "22.77"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("E47").Value = "  +3.48%  "
$ws.Range("E48").Value = "  +3.01%  "
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("E51").Value = "  +0.76%  "
